$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.265.56"
Set-TextValue $ws.Range("E2") "  -0.77%  "
Set-TextValue $ws.Range("D3") "1.702.38"
Set-TextValue $ws.Range("E3") "  -1.28%  "
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "223.44"
Set-TextValue $ws.Range("E5") "  -1.08%  "
Set-TextValue $ws.Range("D6") "0.5307"
Set-TextValue $ws.Range("E6") "  -1.28%  "
Set-TextValue $ws.Range("D7") "1.003"
Set-TextValue $ws.Range("E7") "  -0.08%  "
Set-TextValue $ws.Range("D8") "0.2660"
Set-TextValue $ws.Range("E8") "  -0.84%  "
Set-TextValue $ws.Range("D9") "0.06587"
Set-TextValue $ws.Range("E9") "  -0.15%  "
Set-TextValue $ws.Range("D10") "20.72"
Set-TextValue $ws.Range("E10") "  -4.31%  "
Set-TextValue $ws.Range("D11") "0.07615"
Set-TextValue $ws.Range("E11") "  -1.40%  "
Set-TextValue $ws.Range("D12") "4.493"
Set-TextValue $ws.Range("D13") "1.713.65"
Set-TextValue $ws.Range("E13") "  -0.51%  "
Set-TextValue $ws.Range("D14") "1.936.80"
Set-TextValue $ws.Range("E14") "  -1.21%  "
Set-TextValue $ws.Range("D15") "0.5780"
Set-TextValue $ws.Range("E15") "  -1.54%  "
Set-TextValue $ws.Range("E16") "  -1.60%  "
Set-TextValue $ws.Range("D17") "67.46"
Set-TextValue $ws.Range("E17") "  -0.79%  "
Set-TextValue $ws.Range("D18") "27.273.96"
Set-TextValue $ws.Range("E18") "  -0.81%  "
Set-TextValue $ws.Range("E19") "  -2.83%  "
Set-TextValue $ws.Range("E20") "  -0.10%  "
Set-TextValue $ws.Range("D21") "4.608"
Set-TextValue $ws.Range("E21") "  -2.67%  "
Set-TextValue $ws.Range("E22") "  -3.03%  "
Set-TextValue $ws.Range("D23") "5.964"
Set-TextValue $ws.Range("E23") "  -2.14%  "
Set-TextValue $ws.Range("D24") "1.004"
Set-TextValue $ws.Range("E24") "  -0.10%  "
Set-TextValue $ws.Range("D25") "143.95"
Set-TextValue $ws.Range("E25") "  -2.50%  "
Set-TextValue $ws.Range("D26") "1.703"
Set-TextValue $ws.Range("E26") "  +0.80%  "
Set-TextValue $ws.Range("D27") "0.1199"
Set-TextValue $ws.Range("E27") "  -2.69%  "
Set-TextValue $ws.Range("D28") "7.200"
Set-TextValue $ws.Range("E28") "  -2.81%  "
Set-TextValue $ws.Range("E29") "  -3.33%  "
Set-TextValue $ws.Range("D30") "0.05364"
Set-TextValue $ws.Range("E30") "  -3.35%  "
Set-TextValue $ws.Range("D31") "1.284"
Set-TextValue $ws.Range("D32") "3.458"
Set-TextValue $ws.Range("E32") "  -2.18%  "
Set-TextValue $ws.Range("D33") "3.399"
Set-TextValue $ws.Range("E33") "  -1.82%  "
Set-TextValue $ws.Range("D34") "1.640"
Set-TextValue $ws.Range("E34") "  -1.01%  "
Set-TextValue $ws.Range("D35") "2.863"
Set-TextValue $ws.Range("E35") "  +1.80%  "
Set-TextValue $ws.Range("D36") "2.415"
Set-TextValue $ws.Range("E36") "  -1.42%  "
Set-TextValue $ws.Range("D37") "0.9445"
Set-TextValue $ws.Range("E37") "  -1.50%  "
Set-TextValue $ws.Range("D38") "0.5804"
Set-TextValue $ws.Range("E38") "  -1.87%  "
Set-TextValue $ws.Range("E39") "  -0.91%  "
Set-TextValue $ws.Range("E40") "  -1.44%  "
Set-TextValue $ws.Range("D41") "1.003"
Set-TextValue $ws.Range("E41") "  -0.06%  "
Set-TextValue $ws.Range("D42") "1.040.29"
Set-TextValue $ws.Range("E42") "  -1.32%  "
Set-TextValue $ws.Range("D43") "0.8396"
Set-TextValue $ws.Range("E43") "  -1.88%  "
Set-TextValue $ws.Range("D44") "100.94"
Set-TextValue $ws.Range("E44") "  -0.72%  "
Set-TextValue $ws.Range("D45") "1.844.82"
Set-TextValue $ws.Range("E45") "  -1.16%  "
Set-TextValue $ws.Range("E46") "  -0.32%  "
Set-TextValue $ws.Range("D47") "57.67"
Set-TextValue $ws.Range("E47") "  -2.03%  "
Set-TextValue $ws.Range("D48") "0.4517"
Set-TextValue $ws.Range("E48") "  +1.72%  "
Set-TextValue $ws.Range("D49") "1.005"
Set-TextValue $ws.Range("E49") "  +0.54%  "
Set-TextValue $ws.Range("D50") "8.044"
Set-TextValue $ws.Range("E50") "  -1.87%  "
Set-TextValue $ws.Range("D51") "0.05225"
Set-TextValue $ws.Range("E51") "  -0.96%  "
